$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Name column (C) for several rows with new names
$ws.Range("C4").Value = "Fairaaz"
$ws.Range("C5").Value = "Siddharth"
$ws.Range("C7").Value = "Atharva"
$ws.Range("C8").Value = "Uday"
$ws.Range("C9").Value = "Pratham"

# Update the selected / active cell
$ws.Range("L9").Select()
